$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the max gestational-age values for DELIV, LB and SB rows from 301 to 308
$ws.Range("B2").Value = 308
$ws.Range("B5").Value = 308
$ws.Range("B6").Value = 308

# Clear the explicit formatting that used to be applied to the data rows,
# reverting them to the default "Normal" style
$ws.Range("A2:D7").Style = "Normal"
